$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enable iterative calculation (used e.g. for the new 2FA cancel-button
# circular-reference guard) with a smaller max change tolerance.
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# Company Code values updated (12323 -> 123236) for the sub-constructor
# account rows. Typed with a leading apostrophe so the cells keep being
# stored/forced as text (matches the existing quote-prefixed formatting).
$ws.Range("D2").Value = "'123236"
$ws.Range("D3").Value = "'123236"
$ws.Range("D4").Value = "'123236"

# Move the saved selection/active cell like the workbook was left in the UI.
$ws.Range("H12").Select()
